$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.209896683692932
$ws.Range("B1").Value = 1.487979173660278
$ws.Range("C1").Value = 6.946109294891357
$ws.Range("D1").Value = 2.191200971603394
$ws.Range("E1").Value = 1.171060562133789
